$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colB = New-Object 'object[,]' 24,1
$colB[0,0] = 0.1021555616740102
$colB[1,0] = 0.09096027898117143
$colB[2,0] = 0.08408673090825403
$colB[3,0] = 0.08128605533616451
$colB[4,0] = 0.08082103386608708
$colB[5,0] = 0.08404895826434711
$colB[6,0] = 0.09829549157478823
$colB[7,0] = 0.1262264331171252
$colB[8,0] = 0.1467323726724317
$colB[9,0] = 0.1560556089372369
$colB[10,0] = 0.1595851251667568
$colB[11,0] = 0.1588250290536024
$colB[12,0] = 0.156346005763865
$colB[13,0] = 0.1548273953757899
$colB[14,0] = 0.1461229532344817
$colB[15,0] = 0.1407815827380858
$colB[16,0] = 0.1377089114036494
$colB[17,0] = 0.1366684875325461
$colB[18,0] = 0.141350229923944
$colB[19,0] = 0.1570741840356789
$colB[20,0] = 0.1673447960835119
$colB[21,0] = 0.1618638059317732
$colB[22,0] = 0.1410931502617814
$colB[23,0] = 0.1186722030635963
$ws.Range("B2:B25").Value = $colB

$colD = New-Object 'object[,]' 24,1
$colD[0,0] = 0.01642777940779894
$colD[1,0] = 0.01441394226924331
$colD[2,0] = 0.01317131093117752
$colD[3,0] = 0.01266342406631082
$colD[4,0] = 0.01257900003751189
$colD[5,0] = 0.01316446743543764
$colD[6,0] = 0.01573470119191711
$colD[7,0] = 0.02072494182259987
$colD[8,0] = 0.02435930573314948
$colD[9,0] = 0.02600543359064034
$colD[10,0] = 0.02662771676664022
$colD[11,0] = 0.02649374515884517
$colD[12,0] = 0.02605665086861819
$colD[13,0] = 0.02578877750683262
$colD[14,0] = 0.02425157987875082
$colD[15,0] = 0.02330669612549485
$colD[16,0] = 0.02276255218527012
$colD[17,0] = 0.02257820022710177
$colD[18,0] = 0.02340735043005537
$colD[19,0] = 0.02618506540572696
$colD[20,0] = 0.02799420916196738
$colD[21,0] = 0.0270292204146827
$colD[22,0] = 0.0233618474750017
$colD[23,0] = 0.01938045516200759
$ws.Range("D2:D25").Value = $colD

$colE = New-Object 'object[,]' 24,1
$colE[0,0] = 0.4249720752661972
$colE[1,0] = 0.3707566556344091
$colE[2,0] = 0.3375480711026739
$colE[3,0] = 0.3240338622339891
$colE[4,0] = 0.3217909143882167
$colE[5,0] = 0.3373657406246622
$colE[6,0] = 0.4062610755421474
$colE[7,0] = 0.5420740011794152
$colE[8,0] = 0.6424015897960373
$colE[9,0] = 0.6881867887667141
$colE[10,0] = 0.7055472445210711
$colE[11,0] = 0.7018073285057795
$colE[12,0] = 0.6896145816037205
$colE[13,0] = 0.682149163633369
$colE[14,0] = 0.6394124912066843
$colE[15,0] = 0.6132334056088098
$colE[16,0] = 0.5981895005315465
$colE[17,0] = 0.5930981807455993
$colE[18,0] = 0.6160187932416363
$colE[19,0] = 0.6931952605309988
$colE[20,0] = 0.7437670149723488
$colE[21,0] = 0.7167632493378733
$colE[22,0] = 0.6147594984657587
$colE[23,0] = 0.5052453835082531
$ws.Range("E2:E25").Value = $colE

$colF = New-Object 'object[,]' 24,1
$colF[0,0] = 0.4036106788199163
$colF[1,0] = 0.3976163764266119
$colF[2,0] = 0.3942437459020098
$colF[3,0] = 0.3929466644707063
$colF[4,0] = 0.392735949600123
$colF[5,0] = 0.3942259402257022
$colF[6,0] = 0.4014798583325714
$colF[7,0] = 0.4181550460171337
$colF[8,0] = 0.4319121419258067
$colF[9,0] = 0.4385003018318514
$colF[10,0] = 0.4410426963418246
$colF[11,0] = 0.4404930267342948
$colF[12,0] = 0.4387085113012148
$colF[13,0] = 0.43762164760318
$colF[14,0] = 0.4314882431650489
$colF[15,0] = 0.4278102224554914
$colF[16,0] = 0.4257257733269952
$colF[17,0] = 0.4250253431692457
$colF[18,0] = 0.4281985392815812
$colF[19,0] = 0.4392313736643843
$colF[20,0] = 0.4467194949635171
$colF[21,0] = 0.4426975006955871
$colF[22,0] = 0.4280228875133432
$colF[23,0] = 0.4133803978028041
$ws.Range("F2:F25").Value = $colF

$colG = New-Object 'object[,]' 24,1
$colG[0,0] = 0.2551188872457146
$colG[1,0] = 0.2503476832898599
$colG[2,0] = 0.2476528939344433
$colG[3,0] = 0.2466135479801963
$colG[4,0] = 0.2464445102718571
$colG[5,0] = 0.2476386391638385
$colG[6,0] = 0.253424915056371
$colG[7,0] = 0.2666456529680801
$colG[8,0] = 0.2775182862675365
$colG[9,0] = 0.2827200158037186
$colG[10,0] = 0.2847268193103218
$colG[11,0] = 0.2842929677488115
$colG[12,0] = 0.2828843733391011
$colG[13,0] = 0.2820263972947998
$colG[14,0] = 0.277183508004299
$colG[15,0] = 0.2742782247136972
$colG[16,0] = 0.272631236676844
$colG[17,0] = 0.2720777185196823
$colG[18,0] = 0.2745850057614376
$colG[19,0] = 0.2832971049183755
$colG[20,0] = 0.2892068706054687
$colG[21,0] = 0.2860328794648979
$colG[22,0] = 0.2744462374953969
$colG[23,0] = 0.2628667260298698
$ws.Range("G2:G25").Value = $colG

$colH = New-Object 'object[,]' 24,1
$colH[0,0] = 0.4115050326584253
$colH[1,0] = 0.412803309531192
$colH[2,0] = 0.413826878323988
$colH[3,0] = 0.4143009011998302
$colH[4,0] = 0.4143830488538853
$colH[5,0] = 0.413833040769461
$colH[6,0] = 0.4119056686234757
$colH[7,0] = 0.4099241975833934
$colH[8,0] = 0.4095673111052918
$colH[9,0] = 0.4096442091024244
$colH[10,0] = 0.4097077757199941
$colH[11,0] = 0.4096925527716877
$colH[12,0] = 0.4096487481456421
$colH[13,0] = 0.4096264038665964
$colH[14,0] = 0.4095671036183859
$colH[15,0] = 0.4095920326435589
$colH[16,0] = 0.4096288876187515
$colH[17,0] = 0.4096452319882928
$colH[18,0] = 0.4095870482219368
$colH[19,0] = 0.4096606793940794
$colH[20,0] = 0.4099096057802427
$colH[21,0] = 0.4097583609674018
$colH[22,0] = 0.4095892315217071
$colH[23,0] = 0.4102674590303934
$ws.Range("H2:H25").Value = $colH

$colI = New-Object 'object[,]' 24,1
$colI[0,0] = 0.347860017465653
$colI[1,0] = 0.3532795260296455
$colI[2,0] = 0.3567975250746325
$colI[3,0] = 0.3582790100338156
$colI[4,0] = 0.3585279006596005
$colI[5,0] = 0.35681731107494
$colI[6,0] = 0.3496891439856389
$colI[7,0] = 0.3372216081716659
$colI[8,0] = 0.3289824973184466
$colI[9,0] = 0.3254341973865075
$colI[10,0] = 0.3241192730929507
$colI[11,0] = 0.3244011870107303
$colI[12,0] = 0.3253254416603873
$colI[13,0] = 0.3258953177389423
$colI[14,0] = 0.3292184076257653
$colI[15,0] = 0.3313081821024602
$colI[16,0] = 0.3325289617470004
$colI[17,0] = 0.3329455244493759
$colI[18,0] = 0.3310837765055822
$colI[19,0] = 0.325053185502228
$colI[20,0] = 0.321279373705041
$colI[21,0] = 0.3232781921560158
$colI[22,0] = 0.3311851700335063
$colI[23,0] = 0.3404326915049554
$ws.Range("I2:I25").Value = $colI

$colK = New-Object 'object[,]' 24,1
$colK[0,0] = 0.6551946517829492
$colK[1,0] = 0.5740045889344572
$colK[2,0] = 0.5239262514584482
$colK[3,0] = 0.5034627838914787
$colK[4,0] = 0.5000614801823531
$colK[5,0] = 0.5236504999250826
$colK[6,0] = 0.6272480843206267
$colK[7,0] = 0.8285636016978799
$colK[8,0] = 0.9753152803850469
$colK[9,0] = 1.041819376314862
$colK[10,0] = 1.066965378893485
$colK[11,0] = 1.061551430056284
$colK[12,0] = 1.04388891422829
$colK[13,0] = 1.033065175884389
$colK[14,0] = 0.9709638736522095
$colK[15,0] = 0.9328008787005331
$colK[16,0] = 0.9108266952485167
$colK[17,0] = 0.9033825546642333
$colK[18,0] = 0.9368658700410037
$colK[19,0] = 1.049077854241062
$colK[20,0] = 1.122194865550512
$colK[21,0] = 1.083191444861654
$colK[22,0] = 0.9350281925621289
$colK[23,0] = 0.7743027973165226
$ws.Range("K2:K25").Value = $colK

$colN = New-Object 'object[,]' 24,1
$colN[0,0] = 1.02715882567027
$colN[1,0] = 1.033870049886879
$colN[2,0] = 1.038371686572646
$colN[3,0] = 1.040302059196634
$colN[4,0] = 1.040628393382264
$colN[5,0] = 1.03839733170124
$colN[6,0] = 1.029393875775455
$colN[7,0] = 1.014754879565679
$colN[8,0] = 1.005831442227979
$colN[9,0] = 1.002168251796853
$colN[10,0] = 1.00083794742541
$colN[11,0] = 1.001121924623078
$colN[12,0] = 1.00205766770798
$colN[13,0] = 1.0026382402474
$colN[14,0] = 1.0060788031154
$colN[15,0] = 1.008290864941578
$colN[16,0] = 1.009600475174864
$colN[17,0] = 1.01005029440929
$colN[18,0] = 1.008051528433086
$colN[19,0] = 1.001781274694217
$colN[20,0] = 0.9980147197039457
$colN[21,0] = 0.9999947070929451
$colN[22,0] = 1.008159614507115
$colN[23,0] = 1.018392915543963
$ws.Range("N2:N25").Value = $colN

$colO = New-Object 'object[,]' 24,1
$colO[0,0] = 1.251641335899564
$colO[1,0] = 1.243984724533775
$colO[2,0] = 1.240246089720756
$colO[3,0] = 1.238964067453466
$colO[4,0] = 1.238765758968825
$colO[5,0] = 1.240227822810098
$colO[6,0] = 1.24880122582303
$colO[7,0] = 1.27327845085577
$colO[8,0] = 1.295976703546586
$colO[9,0] = 1.307335870868286
$colO[10,0] = 1.311786601228278
$colO[11,0] = 1.31082140842193
$colO[12,0] = 1.307699040140278
$colO[13,0] = 1.305805956130456
$colO[14,0] = 1.295255207876096
$colO[15,0] = 1.289047845896221
$colO[16,0] = 1.285574786928436
$colO[17,0] = 1.284415552787578
$colO[18,0] = 1.289698560468821
$colO[19,0] = 1.308612099507769
$colO[20,0] = 1.321843485336103
$colO[21,0] = 1.314701806156137
$colO[22,0] = 1.289404074548571
$colO[23,0] = 1.265831789532569
$ws.Range("O2:O25").Value = $colO

